$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename H1 header
$ws.Range("H1").Value = "Avg_Experiment_Time"

# Insert 6 new columns before the old "Obs_Prob" column (I), pushing
# Obs_Prob (and the old, now-unused Std_Total_Rounds column after it)
# to the right.
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()

# The old "Std_Total_Rounds" column (previously J) is no longer used in
# this position - it has been replaced by the richer set of Std_* columns
# below, so drop the now-orphaned column (currently shifted to P).
$ws.Range("P1").EntireColumn.Delete()

# New column headers
$ws.Range("I1").Value = "Std_Total_Rounds"
$ws.Range("J1").Value = "Std_Expl_Cost"
$ws.Range("K1").Value = "Std_Expl_Eff"
$ws.Range("L1").Value = "Std_Round_Time"
$ws.Range("M1").Value = "Std_Agent_Step_Time"
$ws.Range("N1").Value = "Std_Experiment_Time"

# Data rows 2-13 (rows 2/3, 6/7 and 10/11 are swapped relative to the
# original sheet, values refreshed with full precision, and the new
# Std_* columns populated).
# row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 52.474
$ws.Range("D2").Value = 52.474
$ws.Range("E2").Value = 3.27932476
$ws.Range("F2").Value = 0.11982008
$ws.Range("G2").Value = 0.11982008
$ws.Range("H2").Value = 6.28757532
$ws.Range("I2").Value = 6.93639323439815
$ws.Range("J2").Value = 6.93639323439815
$ws.Range("K2").Value = 0.4483329319748089
$ws.Range("L2").Value = 0.0005816588757127249
$ws.Range("M2").Value = 0.0005816588757127249
$ws.Range("N2").Value = 0.8325617832983421
$ws.Range("O2").Value = 0.15
# row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 88.318
$ws.Range("D3").Value = 88.318
$ws.Range("E3").Value = 1.95153118
$ws.Range("F3").Value = 0.1234054
$ws.Range("G3").Value = 0.1234054
$ws.Range("H3").Value = 10.89602058
$ws.Range("I3").Value = 12.41406839792567
$ws.Range("J3").Value = 12.41406839792567
$ws.Range("K3").Value = 0.2764158466060103
$ws.Range("L3").Value = 0.004657698439089739
$ws.Range("M3").Value = 0.004657698439089739
$ws.Range("N3").Value = 1.567666006198464
$ws.Range("O3").Value = 0.85
# row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 28.758
$ws.Range("D4").Value = 57.516
$ws.Range("E4").Value = 3.02658988
$ws.Range("F4").Value = 0.12946328
$ws.Range("G4").Value = 0.06473174
$ws.Range("H4").Value = 1.86037352
$ws.Range("I4").Value = 5.126010922340285
$ws.Range("J4").Value = 10.25202184468057
$ws.Range("K4").Value = 0.5080511998826398
$ws.Range("L4").Value = 0.006538954924497432
$ws.Range("M4").Value = 0.003269520561459067
$ws.Range("N4").Value = 0.3357931835309068
$ws.Range("O4").Value = 0.15
# row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 48.394
$ws.Range("D5").Value = 96.786
$ws.Range("E5").Value = 1.80240734
$ws.Range("F5").Value = 0.12183178
$ws.Range("G5").Value = 0.06091596000000001
$ws.Range("H5").Value = 2.94689314
$ws.Range("I5").Value = 8.657058005270223
$ws.Range("J5").Value = 17.30801371265214
$ws.Range("K5").Value = 0.3236011440591599
$ws.Range("L5").Value = 0.001622370365485518
$ws.Range("M5").Value = 0.0008108836803673955
$ws.Range("N5").Value = 0.5221501234368816
$ws.Range("O5").Value = 0.85
# row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 15.366
$ws.Range("D6").Value = 61.464
$ws.Range("E6").Value = 2.92074236
$ws.Range("F6").Value = 0.12576102
$ws.Range("G6").Value = 0.03144018
$ws.Range("H6").Value = 0.48307418
$ws.Range("I6").Value = 4.040382826759557
$ws.Range("J6").Value = 16.16153130703823
$ws.Range("K6").Value = 0.6841049143086116
$ws.Range("L6").Value = 0.002160795324056615
$ws.Range("M6").Value = 0.0005401285451543745
$ws.Range("N6").Value = 0.1272114305180021
$ws.Range("O6").Value = 0.15
# row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 25.416
$ws.Range("D7").Value = 101.644
$ws.Range("E7").Value = 1.74044188
$ws.Range("F7").Value = 0.12768986
$ws.Range("G7").Value = 0.03192232
$ws.Range("H7").Value = 0.8105630000000001
$ws.Range("I7").Value = 5.691020948729476
$ws.Range("J7").Value = 22.75928897107592
$ws.Range("K7").Value = 0.3655158606644315
$ws.Range("L7").Value = 0.007448395624236826
$ws.Range("M7").Value = 0.001861791705403324
$ws.Range("N7").Value = 0.1837964228288717
$ws.Range("O7").Value = 0.85
# row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 10.032
$ws.Range("D8").Value = 60.192
$ws.Range("E8").Value = 2.96112208
$ws.Range("F8").Value = 0.13137902
$ws.Range("G8").Value = 0.02189646
$ws.Range("H8").Value = 0.21953618
$ws.Range("I8").Value = 2.427913822824468
$ws.Range("J8").Value = 14.56748293694681
$ws.Range("K8").Value = 0.6746925128114
$ws.Range("L8").Value = 0.01104355393007219
$ws.Range("M8").Value = 0.001840451518509512
$ws.Range("N8").Value = 0.05604423960759918
$ws.Range("O8").Value = 0.15
# row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 17.752
$ws.Range("D9").Value = 106.49
$ws.Range("E9").Value = 1.6806774
$ws.Range("F9").Value = 0.12690168
$ws.Range("G9").Value = 0.02115014
$ws.Range("H9").Value = 0.3751084
$ws.Range("I9").Value = 4.447702791356509
$ws.Range("J9").Value = 26.67502603612705
$ws.Range("K9").Value = 0.3946577504908334
$ws.Range("L9").Value = 0.00173120852108779
$ws.Range("M9").Value = 0.0002885194401330166
$ws.Range("N9").Value = 0.09244985035954902
$ws.Range("O9").Value = 0.85
# row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 7.45
$ws.Range("D10").Value = 59.6
$ws.Range("E10").Value = 2.991938779999999
$ws.Range("F10").Value = 0.12238388
$ws.Range("G10").Value = 0.01529804
$ws.Range("H10").Value = 0.11394494
$ws.Range("I10").Value = 1.769112636051779
$ws.Range("J10").Value = 14.15290108841424
$ws.Range("K10").Value = 0.6941270994108966
$ws.Range("L10").Value = 0.001447533091055053
$ws.Range("M10").Value = 0.0001805895444568546
$ws.Range("N10").Value = 0.02699184541364425
$ws.Range("O10").Value = 0.15
# row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 13.386
$ws.Range("D11").Value = 107.042
$ws.Range("E11").Value = 1.68667464
$ws.Range("F11").Value = 0.12104198
$ws.Range("G11").Value = 0.01513032
$ws.Range("H11").Value = 0.20237088
$ws.Range("I11").Value = 3.681863146895187
$ws.Range("J11").Value = 29.40328795705288
$ws.Range("K11").Value = 0.4249842414205284
$ws.Range("L11").Value = 0.004901327483631716
$ws.Range("M11").Value = 0.0006125274502866866
$ws.Range("N11").Value = 0.05547959554168124
$ws.Range("O11").Value = 0.85
# row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 6.036
$ws.Range("D12").Value = 60.36
$ws.Range("E12").Value = 3.00514502
$ws.Range("F12").Value = 0.12356824
$ws.Range("G12").Value = 0.01235692
$ws.Range("H12").Value = 0.07458858
$ws.Range("I12").Value = 1.621098368341563
$ws.Range("J12").Value = 16.21098368341563
$ws.Range("K12").Value = 0.8151755741463118
$ws.Range("L12").Value = 0.001508019079949251
$ws.Range("M12").Value = 0.0001509519346452582
$ws.Range("N12").Value = 0.02008402462341021
$ws.Range("O12").Value = 0.15
# row 13
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 11.038
$ws.Range("D13").Value = 110.28
$ws.Range("E13").Value = 1.64855824
$ws.Range("F13").Value = 0.12084636
$ws.Range("G13").Value = 0.01208478
$ws.Range("H13").Value = 0.1331877
$ws.Range("I13").Value = 3.055360282013532
$ws.Range("J13").Value = 30.49566479259445
$ws.Range("K13").Value = 0.4506542722712557
$ws.Range("L13").Value = 0.004213287972769363
$ws.Range("M13").Value = 0.0004212941122993724
$ws.Range("N13").Value = 0.03655539893493963
$ws.Range("O13").Value = 0.85
